# Refresh cryptocurrency price/volume snapshot (automated data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while keeping it a text cell (matches the source feed,
# which always emits inline/shared strings for these columns) and without
# leaving any numeric formatting behind on the cell.
function Set-TextValue([string]$addr, [string]$val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "42.642.41"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "2.527.78"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue "D5" "315.46"
$ws.Range("E5").Value = "  +3.73%  "
Set-TextValue "D6" "94.98"
$ws.Range("E6").Value = "  -3.13%  "
Set-TextValue "D7" "0.580"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  +0.13%  "
Set-TextValue "D9" "0.538"
$ws.Range("E9").Value = "  -1.42%  "
Set-TextValue "D10" "36.18"
$ws.Range("E10").Value = "  -1.40%  "
$ws.Range("E11").Value = "  -1.27%  "
Set-TextValue "D12" "7.59"
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("E13").Value = "  -2.64%  "
$ws.Range("D14").Value = "2.912.16"
Set-TextValue "D15" "15.57"
$ws.Range("E15").Value = "  +2.80%  "
$ws.Range("D16").Value = "2.476.55"
$ws.Range("E16").Value = "  -2.64%  "
Set-TextValue "D17" "0.862"
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").Value = "42.683.97"
$ws.Range("E18").Value = "  -0.61%  "
Set-TextValue "D19" "13.00"
$ws.Range("E19").Value = "  -5.64%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D20" "6.61"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0969"
$ws.Range("E21").Value = "  -2.55%  "
Set-TextValue "D22" "71.35"
$ws.Range("E22").Value = "  -0.79%  "
Set-TextValue "D23" "255.16"
$ws.Range("E23").Value = "  +0.28%  "
Set-TextValue "D24" "2.97"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  -2.05%  "
Set-TextValue "D26" "27.71"
$ws.Range("E26").Value = "  -1.51%  "
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("E28").Value = "  +10.58%  "
Set-TextValue "D29" "39.26"
$ws.Range("E29").Value = "  +3.76%  "
Set-TextValue "D30" "10.06"
$ws.Range("E30").Value = "  -2.03%  "
Set-TextValue "D31" "5.90"
$ws.Range("E31").Value = "  -4.53%  "
Set-TextValue "D32" "155.96"
$ws.Range("E32").Value = "  -1.74%  "
Set-TextValue "D33" "19.98"
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("E34").Value = "  +0.92%  "
Set-TextValue "D35" "2.11"
$ws.Range("E35").Value = "  -2.33%  "
$ws.Range("E36").Value = "  -2.18%  "
Set-TextValue "D37" "2.61"
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("E38").Value = "  -3.37%  "
Set-TextValue "D39" "24.59"
$ws.Range("E39").Value = "  -4.11%  "
Set-TextValue "D40" "0.120"
$ws.Range("E40").Value = "  -0.15%  "
Set-TextValue "D41" "2.19"
$ws.Range("E41").Value = "  +3.26%  "
$ws.Range("E42").Value = "  -1.06%  "
Set-TextValue "D43" "3.85"
$ws.Range("E43").Value = "  -1.54%  "
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").Value = "2.051.29"
$ws.Range("E46").Value = "  -1.91%  "
Set-TextValue "D47" "86.34"
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.768.87"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
Set-TextValue "D50" "74.44"
$ws.Range("E50").Value = "  -1.33%  "
Set-TextValue "D51" "0.189"
$ws.Range("E51").Value = "  -0.50%  "
